$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 6.246
$ws.Range("B6").Value = 6.728
$ws.Range("B7").Value = 6.638
$ws.Range("B8").Value = 5.726
$ws.Range("B16").Value = 6.184
$ws.Range("B20").Value = 5.776
$ws.Range("B21").Value = 6.247
